# Scheduled runner update: refresh market-derived profit columns (H:N) per leve row.
# Generated from the authoritative cell-level diff; one block per worksheet.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H99").Value = 532.875
$ws.Range("J99").Value = 625.5
$ws.Range("L99").Value = 1876.5
$ws.Range("N99").Value = -4872.5
$ws.Range("H101").Value = 1680.375
$ws.Range("I101").Value = 243.25
$ws.Range("J101").Value = 3117.5
$ws.Range("K101").Value = 729.75
$ws.Range("L101").Value = 9352.5
$ws.Range("M101").Value = 892.25
$ws.Range("N101").Value = -12596.5
$ws.Range("H115").Value = 626
$ws.Range("I115").Value = 626
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1878
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -311
$ws.Range("H118").Value = 836.6
$ws.Range("I118").Value = 356
$ws.Range("J118").Value = 1042.5714
$ws.Range("K118").Value = 1068
$ws.Range("L118").Value = 3127.7142
$ws.Range("M118").Value = 589
$ws.Range("N118").Value = -6441.7142
$ws.Range("H123").Value = 58329.23
$ws.Range("J123").Value = 58329.23
$ws.Range("L123").Value = 58329.23
$ws.Range("N123").Value = -68129.23000000001
$ws.Range("H124").Value = 42968.57
$ws.Range("J124").Value = 42968.57
$ws.Range("L124").Value = 42968.57
$ws.Range("N124").Value = -52788.57
$ws.Range("H126").Value = 28864.445
$ws.Range("J126").Value = 28864.445
$ws.Range("L126").Value = 28864.445
$ws.Range("N126").Value = -38744.445
$ws.Range("H127").Value = 1753.591
$ws.Range("I127").Value = 698
$ws.Range("J127").Value = 1988.1666
$ws.Range("K127").Value = 2094
$ws.Range("L127").Value = 5964.4998
$ws.Range("M127").Value = 2866
$ws.Range("N127").Value = -15884.4998
$ws.Range("H128").Value = 29900
$ws.Range("J128").Value = 29900
$ws.Range("L128").Value = 29900
$ws.Range("N128").Value = -39860
$ws.Range("H129").Value = 954.8289
$ws.Range("I129").Value = 394.25
$ws.Range("J129").Value = 985.9722
$ws.Range("K129").Value = 1182.75
$ws.Range("L129").Value = 2957.9166
$ws.Range("M129").Value = 3817.25
$ws.Range("N129").Value = -12957.9166
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("N130").Value = 0
$ws.Range("H131").Value = 3933
$ws.Range("I131").Value = 3933
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 11799
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -6759
$ws.Range("H133").Value = 21571.428
$ws.Range("J133").Value = 21571.428
$ws.Range("L133").Value = 21571.428
$ws.Range("N133").Value = -31691.428
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = 0
$ws.Range("H137").Value = 1941.8788
$ws.Range("I137").Value = 1366.2142
$ws.Range("J137").Value = 2366.0527
$ws.Range("K137").Value = 4098.642599999999
$ws.Range("L137").Value = 7098.158100000001
$ws.Range("M137").Value = -1548.642599999999
$ws.Range("N137").Value = -12198.1581
$ws.Range("H138").Value = 3395.4092
$ws.Range("I138").Value = 1597.091
$ws.Range("J138").Value = 4294.5684
$ws.Range("K138").Value = 4791.272999999999
$ws.Range("L138").Value = 12883.7052
$ws.Range("M138").Value = 348.7270000000008
$ws.Range("N138").Value = -23163.7052
$ws.Range("H140").Value = 34000
$ws.Range("J140").Value = 34000
$ws.Range("L140").Value = 34000
$ws.Range("N140").Value = -44360

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 34694.5
$ws.Range("J52").Value = 34694.5
$ws.Range("L52").Value = 34694.5
$ws.Range("N52").Value = -35330.5
$ws.Range("H122").Value = 1071692.2
$ws.Range("I122").Value = 1285299.2
$ws.Range("J122").Value = 3657
$ws.Range("K122").Value = 3855897.6
$ws.Range("L122").Value = 10971
$ws.Range("M122").Value = -3853447.6
$ws.Range("N122").Value = -15871

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4066.0417
$ws.Range("I134").Value = 5790.8335
$ws.Range("J134").Value = 2341.25
$ws.Range("K134").Value = 17372.5005
$ws.Range("L134").Value = 7023.75
$ws.Range("M134").Value = -14837.5005
$ws.Range("N134").Value = -12093.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2284.0117
$ws.Range("I31").Value = 1612.6316
$ws.Range("J31").Value = 2474.403
$ws.Range("K31").Value = 1612.6316
$ws.Range("L31").Value = 2474.403
$ws.Range("M31").Value = -1317.6316
$ws.Range("N31").Value = -3064.403
$ws.Range("H34").Value = 2284.0117
$ws.Range("I34").Value = 1612.6316
$ws.Range("J34").Value = 2474.403
$ws.Range("K34").Value = 1612.6316
$ws.Range("L34").Value = 2474.403
$ws.Range("M34").Value = -1410.6316
$ws.Range("N34").Value = -2878.403
$ws.Range("H56").Value = 27103
$ws.Range("J56").Value = 27103
$ws.Range("L56").Value = 27103
$ws.Range("N56").Value = -28793
$ws.Range("H132").Value = 2997.3333
$ws.Range("I132").Value = 1918.4615
$ws.Range("J132").Value = 5802.4
$ws.Range("K132").Value = 5755.3845
$ws.Range("L132").Value = 17407.2
$ws.Range("M132").Value = -3225.3845
$ws.Range("N132").Value = -22467.2
$ws.Range("H133").Value = 39975
$ws.Range("J133").Value = 39975
$ws.Range("L133").Value = 39975
$ws.Range("N133").Value = -45035

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 360.30435
$ws.Range("I40").Value = 190.9375
$ws.Range("J40").Value = 747.4286
$ws.Range("K40").Value = 763.75
$ws.Range("L40").Value = 2989.7144
$ws.Range("M40").Value = -694.75
$ws.Range("N40").Value = -3127.7144
$ws.Range("H68").Value = 2004.8649
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2004.8649
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").Value = 6014.5947
$ws.Range("N68").Value = -7636.5947
$ws.Range("H71").Value = 2004.8649
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2004.8649
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").Value = 18043.7841
$ws.Range("N71").Value = -26155.7841
$ws.Range("H107").Value = 1176.46
$ws.Range("I107").Value = 351.75
$ws.Range("J107").Value = 1333.5476
$ws.Range("K107").Value = 1055.25
$ws.Range("L107").Value = 4000.642800000001
$ws.Range("M107").Value = 864.75
$ws.Range("N107").Value = -7840.642800000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("K97").Value = 1000
$ws.Range("M97").Value = -504
$ws.Range("H122").Value = 5402678
$ws.Range("I122").Value = 6482965
$ws.Range("J122").Value = 1245
$ws.Range("K122").Value = 19448895
$ws.Range("L122").Value = 3735
$ws.Range("M122").Value = -19446445
$ws.Range("N122").Value = -8635
$ws.Range("H126").Value = 6597.857
$ws.Range("I126").Value = 8095.3125
$ws.Range("J126").Value = 1806
$ws.Range("K126").Value = 24285.9375
$ws.Range("L126").Value = 5418
$ws.Range("M126").Value = -21815.9375
$ws.Range("N126").Value = -10358

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12347274
$ws.Range("I22").Value = 55556156
$ws.Range("J22").Value = 1878.5714
$ws.Range("K22").Value = 55556156
$ws.Range("L22").Value = 1878.5714
$ws.Range("M22").Value = -55555861
$ws.Range("N22").Value = -2468.5714
$ws.Range("H27").Value = 12347274
$ws.Range("I27").Value = 55556156
$ws.Range("J27").Value = 1878.5714
$ws.Range("K27").Value = 55556156
$ws.Range("L27").Value = 1878.5714
$ws.Range("M27").Value = -55556049
$ws.Range("N27").Value = -2092.5714
$ws.Range("H46").Value = 18519536
$ws.Range("I46").Value = 41667630
$ws.Range("J46").Value = 1060.1
$ws.Range("K46").Value = 41667630
$ws.Range("L46").Value = 1060.1
$ws.Range("M46").Value = -41667442
$ws.Range("N46").Value = -1436.1
$ws.Range("H61").Value = 2683.0833
$ws.Range("I61").Value = 2021.3334
$ws.Range("K61").Value = 2021.3334
$ws.Range("M61").Value = -1819.3334
$ws.Range("H113").Value = 2683.0833
$ws.Range("I113").Value = 2021.3334
$ws.Range("K113").Value = 2021.3334
$ws.Range("M113").Value = 148.6666
$ws.Range("H136").Value = 6163.4907
$ws.Range("I136").Value = 4702
$ws.Range("K136").Value = 14106
$ws.Range("M136").Value = -11556
$ws.Range("H141").Value = 49800
$ws.Range("J141").Value = 49800
$ws.Range("L141").Value = 49800
$ws.Range("N141").Value = -60160

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2354.5715
$ws.Range("I96").Value = 1674.75
$ws.Range("J96").Value = 2626.5
$ws.Range("K96").Value = 1674.75
$ws.Range("L96").Value = 2626.5
$ws.Range("M96").Value = -301.75
$ws.Range("N96").Value = -5372.5
$ws.Range("H113").Value = 972.2759
$ws.Range("I113").Value = 703.04
$ws.Range("J113").Value = 2655
$ws.Range("K113").Value = 2109.12
$ws.Range("L113").Value = 7965
$ws.Range("M113").Value = 60.88000000000011
$ws.Range("N113").Value = -12305
$ws.Range("H122").Value = 5708.7856
$ws.Range("I122").Value = 6163.5
$ws.Range("J122").Value = 5102.5
$ws.Range("K122").Value = 18490.5
$ws.Range("L122").Value = 15307.5
$ws.Range("M122").Value = -16040.5
$ws.Range("N122").Value = -20207.5

